$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the "Direct By Register" label in A8 to "By Register"
$ws.Range("A8").Value = "By Register"

# Add a new narrow column I holding usage examples next to each addressing-mode block.
$ws.Range("I5").Value = "e.g. LD A, [0x00c0]"
$ws.Range("I8").Value = "e.g. LD A, C"
$ws.Range("I11").Value = "e.g. LD H, 0xff"
$ws.Range("I14").Value = "e.g. LD A, [HL]"

# Bold + left-align the example column cells for each addressing-mode block
$exampleRanges = @("I5:I7", "I8:I10", "I11:I13", "I14:I15")
foreach ($rng in $exampleRanges) {
    $r = $ws.Range($rng)
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4131  # xlLeft
}

# Give column I a sensible default width
$ws.Columns("I").ColumnWidth = 8.89

# Widen column H so the longer "e.g." examples (and existing long text) fit.
$ws.Columns("H").ColumnWidth = 31.11

$ws.Range("J17").Select()
